# Add Composite IDs to the tck sheet.
#   - Row 19 (was "reserved") becomes "composite"
#   - Row 20 (was "reserved") becomes "composite-list", with bit flags for
#     the "name" (E, value 8), "int value" (G omitted / not used here) and
#     other columns set so the total-type formula recalculates to 72.
#   - Row 21 (was "reserved") becomes "composite-map", with bit flags set so
#     the total-type formula recalculates to 80.
#   - Row 15's E/F cells are cleared back to blank (no stored 0 value).
#   - The active selection moves from E2 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: clear the stray zero values in E15/F15 so they go back to empty cells.
$ws.Range("E15").ClearContents()
$ws.Range("F15").ClearContents()

# Row 19: rename from "reserved" to "composite" (flag bits unchanged).
$ws.Range("A19").Value = "composite"

# Row 20: rename from "reserved" to "composite-list" and set its flag bits.
$ws.Range("A20").Value = "composite-list"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 1

# Row 21: rename from "reserved" to "composite-map" and set its flag bits.
$ws.Range("A21").Value = "composite-map"
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 1

# Move the sheet's active cell/selection from E2 to B2.
[void]$ws.Range("B2").Select()
